# Add a new "Exceptional items" column to the Quarterly results sheet,
# between the existing "P/l before exceptional items & tax" column (K)
# and "P/l before tax" column (old L, now shifted to M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new column at L — this shifts the old L:T data right to M:U,
# carrying values/formatting with it (matches the new dimension A1:U47).
$ws.Columns("L").Insert()

# Header rows for the newly inserted column.
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"

# Populate the exceptional-items values for the quarters that reported one
# (difference between "P/l before int., excpt. items & tax" and
# "P/l before tax" in the source data).
$ws.Range("L7").Value = -6.44
$ws.Range("L11").Value = -10.44
$ws.Range("L14").Value = -3.98
$ws.Range("L16").Value = -31.63
$ws.Range("L17").Value = -0.82
$ws.Range("L22").Value = -13.08
